$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28; this shifts the existing rows 28-84
# down to 29-85 (preserving all of their data/formatting), matching the
# diff where every row from the former row 28 onward moved down by one
# and a brand-new row appeared at row 28.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted (blank) row 28 with its data.
$ws.Cells.Item(28,1).Value2 = 8
$ws.Cells.Item(28,2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(28,3).Value2 = "Coquimbo"
$ws.Cells.Item(28,4).Value2 = 44469
$ws.Cells.Item(28,5).Value2 = 4
$ws.Cells.Item(28,6).Value2 = 100112044
$ws.Cells.Item(28,7).Value2 = "Perejil"
$ws.Cells.Item(28,8).Value2 = "Sin especificar"
$ws.Cells.Item(28,9).Value2 = "Primera"
$ws.Cells.Item(28,10).Value2 = 3100
$ws.Cells.Item(28,11).Value2 = 1500
$ws.Cells.Item(28,12).Value2 = 2000
$ws.Cells.Item(28,13).Value2 = 1750
$ws.Cells.Item(28,14).Value2 = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(28,15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(28,16).Value2 = 1167
$ws.Cells.Item(28,17).Value2 = 1.5
$ws.Cells.Item(28,18).Value2 = "Hortaliza"
